$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '246.90'
    'E2' = '0.78%'
    'E3' = '5.14%'
    'D4' = '5.072'
    'E4' = '1.29%'
    'E5' = '-0.20%'
    'D6' = '6.473'
    'E6' = '-1.53%'
    'D7' = '0.8134'
    'E7' = '0.33%'
    'D8' = '0.8423'
    'E8' = '0.23%'
    'D9' = '0.06982'
    'E9' = '0.58%'
    'D10' = '0.02864'
    'E10' = '0.87%'
    'D11' = '0.09387'
    'E11' = '-0.21%'
    'D12' = '0.001520'
    'E12' = '-0.38%'
    'D13' = '0.0005980'
    'E13' = '0.21%'
    'D14' = '0.006193'
    'E14' = '1.24%'
    'E15' = '3.05%'
    'E16' = '0.27%'
    'E17' = '-1.73%'
    'E18' = '-2.15%'
    'D19' = '0.1340'
    'E19' = '0.30%'
    'D20' = '0.03192'
    'E20' = '-2.72%'
    'E21' = '0.52%'
    'D22' = '3.740'
    'E22' = '0.02%'
    'D23' = '0.04651'
    'E23' = '-0.46%'
    'E24' = '-1.41%'
    'D25' = '0.001245'
    'E25' = '0.24%'
    'D26' = '0.004587'
    'E26' = '1.41%'
    'D27' = '0.00009599'
    'E27' = '-1.02%'
    'E28' = '-27.92%'
    'D40' = '0.03663'
    'E40' = '0.06%'
    'D41' = '0.1350'
    'E41' = '-0.71%'
    'E42' = '-2.74%'
    'D43' = '0.003422'
    'E43' = '-45.16%'
    'D44' = '0.008272'
    'E44' = '2.74%'
    'D45' = '0.00005394'
    'E45' = '2.29%'
    'E46' = '0.05%'
    'E47' = '-38.86%'
    'D48' = '0.002604'
    'E48' = '27.51%'
    'E49' = '0.05%'
    'E50' = '0.05%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $updates[$cellRef]
    $range.Style = "Normal"
}
